$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.522.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.61%  "

$ws.Range("D3").Value = "'1.725.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.73%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'225.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.53%  "

$ws.Range("D6").Value = "'0.5351"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.16%  "

$ws.Range("D7").Value = "'1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "'0.2669"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.84%  "

$ws.Range("D9").Value = "'0.06596"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.69%  "

$ws.Range("D10").Value = "'21.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.56%  "

$ws.Range("D11").Value = "'0.07713"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.53%  "

$ws.Range("D12").Value = "'4.611"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.76%  "

$ws.Range("D13").Value = "'1.726.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.82%  "

$ws.Range("D14").Value = "'1.963.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.85%  "

$ws.Range("D15").Value = "'0.5826"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.54%  "

$ws.Range("D16").Value = "'0.0₅8286"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.24%  "

$ws.Range("D17").Value = "'67.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.40%  "

$ws.Range("D18").Value = "'27.538.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.77%  "

$ws.Range("D19").Value = "'218.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +13.39%  "

$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("D21").Value = "'4.728"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.97%  "

$ws.Range("D22").Value = "'10.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.97%  "

$ws.Range("D23").Value = "'6.089"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.02%  "

$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("D25").Value = "'144.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.88%  "

$ws.Range("D26").Value = "'1.766"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +16.44%  "

$ws.Range("D27").Value = "'0.1235"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.77%  "

$ws.Range("D28").Value = "'7.402"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.18%  "

$ws.Range("D29").Value = "'16.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.72%  "

$ws.Range("D30").Value = "'0.05526"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.22%  "

$ws.Range("E31").Value = "  +2.91%  "

$ws.Range("D32").Value = "'3.571"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.79%  "

$ws.Range("D33").Value = "'3.442"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.63%  "

$ws.Range("D34").Value = "'1.657"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.19%  "

$ws.Range("D35").Value = "'2.857"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.87%  "

$ws.Range("D36").Value = "'0.9657"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.63%  "

$ws.Range("E37").Value = "  +0.62%  "

$ws.Range("D38").Value = "'0.5984"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.27%  "

$ws.Range("D39").Value = "'0.01649"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.97%  "

$ws.Range("D40").Value = "'5.903"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.25%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "'1.056.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.80%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.8512"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.30%  "

$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("D44").Value = "'101.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.73%  "

$ws.Range("D45").Value = "'1.869.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.79%  "

$ws.Range("D46").Value = "'0.0₈115"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.25%  "

$ws.Range("D47").Value = "'58.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.00%  "

$ws.Range("D48").Value = "'0.4469"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.71%  "

$ws.Range("D49").Value = "'8.231"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.49%  "

$ws.Range("E50").Value = "  +0.15%  "

$ws.Range("D51").Value = "'0.05245"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.85%  "
